$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------
# Header text updates (Volume/Number bump, reporting week date range)
# -------------------------------------------------------------------
$a8 = $ws.Cells.Item(8, 1)
$chars = $a8.Characters(21, 2)
$chars.Text = "17"

$c9 = $ws.Cells.Item(9, 3)
$d1 = $c9.Characters(27, 9)
$d1.Text = "4/21/2025"
$d2 = $c9.Characters(47, 9)
$d2.Text = "4/27/2025"

# -------------------------------------------------------------------
# Row 15 - Rape: F15 becomes the blank-placeholder text cell ("0")
# -------------------------------------------------------------------
$ws.Range("G15").Copy($ws.Range("F15"))

# -------------------------------------------------------------------
# Row 16 - Robbery
# -------------------------------------------------------------------
$ws.Range("G16").Value = 2
$ws.Range("J16").Value = 8
$ws.Range("K16").Value = -50
$ws.Range("M16").Value = -20
$ws.Range("N16").Value = -76.470588235294

# -------------------------------------------------------------------
# Row 17 - Fel. Assault: D17/E17 flip from placeholder text to real
# numbers, rest of the row updates in place.
# -------------------------------------------------------------------
$ws.Range("C17").Copy($ws.Range("D17"))
$ws.Range("H17").Copy($ws.Range("E17"))
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 34
$ws.Range("J17").Value = 37
$ws.Range("K17").Value = -8.108108108108
$ws.Range("L17").Value = 21.428571428571
$ws.Range("M17").Value = 126.666666666667
$ws.Range("N17").Value = 36

# -------------------------------------------------------------------
# Row 18 - Burglary: D18/E18 flip from placeholder text to real
# numbers, rest of the row updates in place.
# -------------------------------------------------------------------
$ws.Range("C18").Copy($ws.Range("D18"))
$ws.Range("H18").Copy($ws.Range("E18"))
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 66.666666666666
$ws.Range("I18").Value = 24
$ws.Range("J18").Value = 8
$ws.Range("K18").Value = 200
$ws.Range("L18").Value = 20
$ws.Range("M18").Value = -36.842105263157
$ws.Range("N18").Value = -71.084337349397

# -------------------------------------------------------------------
# Row 19 - Gr. Larceny
# -------------------------------------------------------------------
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 20
$ws.Range("G19").Value = 20
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 76
$ws.Range("J19").Value = 88
$ws.Range("K19").Value = -13.636363636363
$ws.Range("L19").Value = -17.391304347826
$ws.Range("M19").Value = 55.102040816326
$ws.Range("N19").Value = 40.740740740740

# -------------------------------------------------------------------
# Row 20 - G.L.A.
# -------------------------------------------------------------------
$ws.Range("I20").Value = 8
$ws.Range("K20").Value = -42.857142857142
$ws.Range("L20").Value = -65.217391304347
$ws.Range("M20").Value = -27.272727272727
$ws.Range("N20").Value = -96.347031963470

# -------------------------------------------------------------------
# Row 21 - TOTAL
# -------------------------------------------------------------------
$ws.Range("C21").Value = 10
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = -28.571428571428
$ws.Range("F21").Value = 32
$ws.Range("G21").Value = 34
$ws.Range("H21").Value = -5.882352941176
$ws.Range("I21").Value = 148
$ws.Range("J21").Value = 158
$ws.Range("K21").Value = -6.329113924050
$ws.Range("L21").Value = -12.426035502958
$ws.Range("M21").Value = 23.333333333333
$ws.Range("N21").Value = -63.092269326683

# -------------------------------------------------------------------
# Row 24 - Petit Larceny
# -------------------------------------------------------------------
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 7
$ws.Range("E24").Value = 28.571428571428
$ws.Range("F24").Value = 36
$ws.Range("G24").Value = 38
$ws.Range("H24").Value = -5.263157894736
$ws.Range("I24").Value = 135
$ws.Range("J24").Value = 137
$ws.Range("K24").Value = -1.459854014598
$ws.Range("L24").Value = -14.556962025316
$ws.Range("M24").Value = -15.625

# -------------------------------------------------------------------
# Row 25 - Retail Theft
# -------------------------------------------------------------------
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -50
$ws.Range("G25").Value = 23
$ws.Range("H25").Value = -13.043478260869
$ws.Range("I25").Value = 70
$ws.Range("J25").Value = 73
$ws.Range("K25").Value = -4.109589041095
$ws.Range("L25").Value = -19.540229885057

# -------------------------------------------------------------------
# Row 26 - Misd. Assault
# -------------------------------------------------------------------
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = -25
$ws.Range("F26").Value = 15
$ws.Range("G26").Value = 15
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 70
$ws.Range("J26").Value = 47
$ws.Range("K26").Value = 48.936170212766
$ws.Range("L26").Value = 11.111111111111
$ws.Range("M26").Value = 14.754098360655

# -------------------------------------------------------------------
# Row 27 - UCR Rape*: F27 becomes the blank-placeholder text cell ("0")
# -------------------------------------------------------------------
$ws.Range("G27").Copy($ws.Range("F27"))

# -------------------------------------------------------------------
# Row 28 - Other Sex Crimes: G28/H28 become blank-placeholder text
# -------------------------------------------------------------------
$ws.Range("G27").Copy($ws.Range("G28"))
$ws.Range("H27").Copy($ws.Range("H28"))

# -------------------------------------------------------------------
# Row 31 - Hate Crimes: G31/H31 become blank-placeholder text
# -------------------------------------------------------------------
$ws.Range("G27").Copy($ws.Range("G31"))
$ws.Range("H27").Copy($ws.Range("H31"))
